$wb = $excel.ActiveWorkbook

# --- 1. "EC1 - 15 years" sheet: update I11 and selection ---
$ws2 = $wb.Worksheets.Item("EC1 - 15 years")
$ws2.Range("I11").Value = 1/(24*7)

# --- 2. "EC1 - 15 yr. w. Prev. Maint. 2y" sheet: update I11 (literal) and I13 (formula) ---
$ws3 = $wb.Worksheets.Item("EC1 - 15 yr. w. Prev. Maint. 2y")
$ws3.Range("I11").Value = 1/(24*7)
$ws3.Range("I13").Formula = "=1/(8760)"

# --- 3. Duplicate the 2y sheet to create the new 1y sheet, placed right after it ---
$ws3.Copy($null, $ws3)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "EC1 - 15 yr. w. Prev. Maint. 1y"

# --- 4. Update selections on each sheet, finishing with the new sheet active ---
$ws2.Activate()
$ws2.Range("I11").Select()

$ws3.Activate()
$ws3.Range("I11").Select()

$ws4.Activate()
$ws4.Range("I14").Select()
